# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets carry the same event list, so the same row/value updates apply
# to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 11561
    3  = 11077
    6  = 1003
    8  = 63
    10 = 41
    11 = 10663
    12 = 4115
    13 = 1
    16 = 47
    17 = 40
    18 = 116
    20 = 11105
    21 = 10869
    22 = 11
    23 = 23
    24 = 9
    26 = 23
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
